{"js": "// The document originally has a hyperlinked run with text\n//   https://github.com/jram828/EVIDENCIAS-SENA\n// followed by a trailing empty paragraph. The edit:\n//   1. Turns the hyperlink into plain (non-linked) text, updating the URL\n//      to include the \"/tree/main/JULIAN_ARANGO_AA4_EV03\" suffix.\n//   2. Removes the trailing empty paragraph.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst newText = \"https://github.com/jram828/EVIDENCIAS-SENA/tree/main/JULIAN_ARANGO_AA4_EV03\";\n\n// Step 1: replace the hyperlinked run with a plain text run carrying the\n// new URL. Clearing the range first drops both the hyperlink field and any\n// character formatting (e.g. the \"Hyperlink\" style) that was riding on it,\n// so the inserted text comes back as an ordinary run.\nconst firstParagraph = paragraphs.items[0];\nconst linkRange = firstParagraph.getRange();\nlinkRange.hyperlink = \"\";\nlinkRange.clear();\nawait context.sync();\n\nlinkRange.insertText(newText, Word.InsertLocation.replace);\nawait context.sync();\n\n// Step 2: delete the trailing empty paragraph. Word will not delete the\n// document's very last paragraph mark via Paragraph.delete(), so instead\n// extend a range from the end of the first paragraph to the end of the\n// second (empty) paragraph and delete that range, which merges them and\n// removes the extra paragraph mark.\nparagraphs.load(\"items\");\nawait context.sync();\n\nif (paragraphs.items.length > 1) {\n  const lastParagraph = paragraphs.items[paragraphs.items.length - 1];\n  const fromEndOfFirst = paragraphs.items[0].getRange(\"End\");\n  const toEndOfLast = lastParagraph.getRange(\"End\");\n  const between = fromEndOfFirst.expandTo(toEndOfLast);\n  between.delete();\n  await context.sync();\n}\n", "ps1": "# The document originally has a hyperlinked run with text\n#   https://github.com/jram828/EVIDENCIAS-SENA\n# followed by a trailing empty paragraph. The edit:\n#   1. Turns the hyperlink into plain (non-linked) text, updating the URL\n#      to include the \"/tree/main/JULIAN_ARANGO_AA4_EV03\" suffix.\n#   2. Removes the trailing empty paragraph.\n\n$d = $word.ActiveDocument\n$newUrl = \"https://github.com/jram828/EVIDENCIAS-SENA/tree/main/JULIAN_ARANGO_AA4_EV03\"\n\n# Step 1: replace the hyperlinked run with a plain text run carrying the\n# new URL. Deleting the hyperlink's range removes the hyperlink field (and\n# the character formatting/style that rode along with it), then retyping\n# the text into that now-empty spot inserts it as an ordinary, unstyled run.\nif ($d.Hyperlinks.Count -gt 0) {\n    $link = $d.Hyperlinks.Item(1)\n    $linkRange = $link.Range\n    $linkRange.Delete()\n    $linkRange.InsertAfter($newUrl)\n} else {\n    $firstParagraph = $d.Paragraphs.First\n    $firstParagraph.Range.Text = $newUrl\n}\n\n# Step 2: delete the trailing empty paragraph. Word will not delete the\n# document's very last paragraph mark by calling Delete() on that\n# paragraph's own range, so instead remove the paragraph mark that ends the\n# *first* paragraph, which merges it with the (empty) one after it.\nif ($d.Paragraphs.Count -gt 1) {\n    $firstParagraph = $d.Paragraphs.First\n    $endOfFirst = $firstParagraph.Range.End\n    $markRange = $d.Range($endOfFirst - 1, $endOfFirst)\n    $markRange.Delete()\n}\n"}
